# Generate Report for handoff
#
# The "711514bb-fefd-4675-bc73-e008cd09d5ab.md" / "Handoff transform failed"
# row is removed from every sheet (its hand-off attempt is gone from the
# report), the in-flight hand-off row now points at a new source file
# ("aa77f244-b95d-42cd-9ac2-cf09660a2f91.md") with a freshly generated
# target/xlf artifact and timestamp, and the ".localization-config" row
# shifts up to take the old row's place.

$wb = $excel.ActiveWorkbook

$oldGuid = "0508241f-403d-4f18-9247-f5fe9e374d5f"
$newGuid = "aa77f244-b95d-42cd-9ac2-cf09660a2f91"
$removedGuid = "711514bb-fefd-4675-bc73-e008cd09d5ab"

$oldHash = "46d937e9d7caef2a89bb212c306488fe8ed2ff9d"
$newHash = "b17dc73c46c65899955740cf9a65663cdcc93e89"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/d8552681a4e67d540e8ec441f88575f0a0aea048/e2e"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d8552681a4e67d540e8ec441f88575f0a0aea048/.localization-config"

$newMdDisplay = "$newGuid.md"
$newMdUrl = "$mdBase/$newMdDisplay"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value() = $newMdDisplay

# Drop the "Handoff transform failed" row entirely - everything below
# shifts up (styles travel with the cells automatically).
$ws.Rows("3").Delete()

# Hyperlinks don't renumber themselves when rows move, so rebuild them
# from scratch against the final layout.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMdDisplay) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$newXlfDisplay = "$newGuid.$newHash.zh-cn.xlf"
$newXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d48280182b3176dea094e5881cbcf892f176e56/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$newXlfDisplay"

$ws.Range("A2").Value() = $newMdDisplay
$ws.Range("C2").Value() = $newXlfDisplay
$ws.Range("D2").Value() = "2016-01-18 06:55:09"

$ws.Rows("3").Delete()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMdDisplay) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $newXlfUrl, "", "", $newXlfDisplay) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$newXlfDisplayDe = "$newGuid.$newHash.de-de.xlf"
$newXlfUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89a45bc423b03372fe25692cf2e2e186421b0988/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$newXlfDisplayDe"

$ws.Range("A2").Value() = $newMdDisplay
$ws.Range("C2").Value() = $newXlfDisplayDe
$ws.Range("D2").Value() = "2016-01-18 06:55:19"

$ws.Rows("3").Delete()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMdDisplay) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $newXlfUrlDe, "", "", $newXlfDisplayDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null
